# "Running suites a and b"
# Update the test-run result for TestCase_B1 (row 2) in the "Results"
# column of the "Test Cases" sheet: it was reported as FAIL, but after
# actually running suites a and b it is now reported as SKIP.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Range("D2").Value = "SKIP"
